# Adding Area / Atotal columns to the discharge worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# Area formulas (column G), rows 2-15
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
for ($r = 4; $r -le 15; $r++) {
    $prev = $r - 1
    $ws.Range("G$r").Formula = "=(D$r-D$prev)*B$r/100"
}

# Totals
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# Selection + view tweaks to match the target state
$ws.Range("J2:K2").Select()

$wb.Save()
